# "Generate Report for Handback"
#
# This mirrors the localization tool re-running its report generator after a
# handback: the zh-cn / de-de sheets gain a populated "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" for both rows, the
# Overview sheet's per-language status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and a couple of columns widen to fit the
# newly-populated long filenames.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: status text for both languages / both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item(1)
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E (zh-cn) & F (de-de) widen to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Helper data shared by the zh-cn (sheet 2) and de-de (sheet 3) tabs
# ---------------------------------------------------------------------
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35d7fef48d38ebf883f8de45644733690802c698/e2e/"
$file1 = "1cb304de-56ac-41d9-a991-05aed199b762.md"
$file2 = "5c1c3f67-e5ad-4f48-bf92-827bfdc9a5ac.md"
$hyperlinkColor = 15570276   # matches the existing custom HyperLink style (FF6495ED)

# ---------------------------------------------------------------------
# zh-cn sheet (sheet 2)
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item(2)

$zh.Range("I2").Value = $file1
$zh.Range("J2").Value = "1cb304de-56ac-41d9-a991-05aed199b762.a65dd7f16e30efec1804e8689c83b800917b6b86.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-26 19:04:34"

$zh.Range("I3").Value = $file2
$zh.Range("J3").Value = "5c1c3f67-e5ad-4f48-bf92-827bfdc9a5ac.2f6f1611d67934dd37b2e2a44a97d42d665adc06.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-26 19:04:34"

$zh.Hyperlinks.Add($zh.Range("I2"), $ghBase + $file1, "", "", $file1)
$zh.Hyperlinks.Add($zh.Range("I3"), $ghBase + $file2, "", "", $file2)

$zh.Range("I2").Font.Color = $hyperlinkColor
$zh.Range("I2").Font.Underline = 2
$zh.Range("I3").Font.Color = $hyperlinkColor
$zh.Range("I3").Font.Underline = 2

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet (sheet 3)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item(3)

$de.Range("I2").Value = $file1
$de.Range("J2").Value = "1cb304de-56ac-41d9-a991-05aed199b762.a65dd7f16e30efec1804e8689c83b800917b6b86.de-de.xlf"
$de.Range("K2").Value = "2016-08-26 19:04:40"

$de.Range("I3").Value = $file2
$de.Range("J3").Value = "5c1c3f67-e5ad-4f48-bf92-827bfdc9a5ac.2f6f1611d67934dd37b2e2a44a97d42d665adc06.de-de.xlf"
$de.Range("K3").Value = "2016-08-26 19:04:40"

$de.Hyperlinks.Add($de.Range("I2"), $ghBase + $file1, "", "", $file1)
$de.Hyperlinks.Add($de.Range("I3"), $ghBase + $file2, "", "", $file2)

$de.Range("I2").Font.Color = $hyperlinkColor
$de.Range("I2").Font.Underline = 2
$de.Range("I3").Font.Color = $hyperlinkColor
$de.Range("I3").Font.Underline = 2

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

"Report regenerated for handback"
